$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear inline-string contents for F3, F4, F5 (EFT rows)
$ws.Range("F3").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("F5").ClearContents()

# Clear inline-string contents for F8, F9, F10 (HAVALE rows)
$ws.Range("F8").ClearContents()
$ws.Range("F9").ClearContents()
$ws.Range("F10").ClearContents()

# Row 13 (GELEN SWIFT) updates
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 851,5 TL"
$ws.Range("F13").ClearContents()
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"

# Row 14 (GIDEN SWIFT - Mobil) update
$ws.Range("F14").ClearContents()
